# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 22:18"

# Update country statistics rows (columns B..H)
# Row 4  - Estados Unidos
$ws.Range("B4").Value = 6452833
$ws.Range("C4").Value = 23693
$ws.Range("D4").Value = 3715139
$ws.Range("E4").Value = 2544573
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = 193121

# Row 5 - India
$ws.Range("E5").Value = 883578
$ws.Range("G5").Value = 1008
$ws.Range("H5").Value = 71687

# Row 10 - Sudafrica
$ws.Range("B10").Value = 638517
$ws.Range("C10").Value = 1633
$ws.Range("D10").Value = 563891
$ws.Range("E10").Value = 59737
$ws.Range("G10").Value = 110
$ws.Range("H10").Value = 14889

# Row 64 - Uzbekistan
$ws.Range("B64").Value = 43587
$ws.Range("C64").Value = 294
$ws.Range("D64").Value = 41225
$ws.Range("E64").Value = 2015
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 347

# Row 110 - Republica de Yibuti
$ws.Range("B110").Value = 5388
$ws.Range("C110").Value = 1
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 61

# Row 120 - Ruanda
$ws.Range("B120").Value = 4374
$ws.Range("C120").Value = 25
$ws.Range("D120").Value = 2235
$ws.Range("E120").Value = 2120
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 19

# Row 134 - Angola
$ws.Range("B134").Value = 2965
$ws.Range("C134").Value = 30
$ws.Range("D134").Value = 1198
$ws.Range("E134").Value = 1650

# Row 154 - Togo
$ws.Range("B154").Value = 1488
$ws.Range("C154").Value = 11
$ws.Range("D154").Value = 1106
$ws.Range("E154").Value = 350

# Row 160 - Liberia
$ws.Range("B160").Value = 1307
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 1182
$ws.Range("E160").Value = 43

$wb.Save()
